# Updates cryptos list values (Price / Volume(1h) columns, plus two row
# swaps for Toncoin/PancakeSwap and Aave/InjectiveProtocol) to match the
# latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value parses as a plain number: Excel COM auto-converts
# such strings to numeric values on assignment, so we prefix with a
# leading apostrophe to force text, then reset the style to "Normal" so
# no stray number-format / quote-prefix style is left attached to the cell.
$numericTextCells = @()

$ws.Range("D2").Value = "41.334.54"
$ws.Range("E2").Value = "  +4.75%  "
$ws.Range("D3").Value = "2.237.17"
$ws.Range("E3").Value = "  +3.49%  "
$ws.Range("D4").Value = "'1.00"
$numericTextCells += "D4"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'228.36"
$numericTextCells += "D5"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("D7").Value = "'64.86"
$numericTextCells += "D7"
$ws.Range("E7").Value = "  +0.64%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.407"
$numericTextCells += "D9"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("D10").Value = "'0.0886"
$numericTextCells += "D10"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").Value = "'0.105"
$numericTextCells += "D11"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "2.570.97"
$ws.Range("E12").Value = "  +3.57%  "
$ws.Range("D13").Value = "'16.09"
$numericTextCells += "D13"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "'22.33"
$numericTextCells += "D14"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "'0.828"
$numericTextCells += "D15"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "'5.63"
$numericTextCells += "D16"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "2.239.75"
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("D18").Value = "41.198.18"
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("D19").Value = "'73.88"
$numericTextCells += "D19"
$ws.Range("E19").Value = "  +2.84%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +6.38%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'253.94"
$numericTextCells += "D22"
$ws.Range("E22").Value = "  +9.57%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.39"
$numericTextCells += "D24"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.34"
$numericTextCells += "D25"
$ws.Range("E25").Value = "  -6.86%  "
$ws.Range("D26").Value = "'9.78"
$numericTextCells += "D26"
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("D27").Value = "'172.90"
$numericTextCells += "D27"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("D29").Value = "'20.40"
$numericTextCells += "D29"
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").Value = "'2.85"
$numericTextCells += "D31"
$ws.Range("E31").Value = "  +6.55%  "
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").Value = "'4.68"
$numericTextCells += "D33"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("D34").Value = "'4.85"
$numericTextCells += "D34"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").Value = "'7.21"
$numericTextCells += "D35"
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "'3.82"
$numericTextCells += "D37"
$ws.Range("E37").Value = "  +6.54%  "
$ws.Range("D38").Value = "'2.46"
$numericTextCells += "D38"
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "'0.000234"
$numericTextCells += "D40"
$ws.Range("E40").Value = "  +51.94%  "
$ws.Range("D41").Value = "'4.81"
$numericTextCells += "D41"
$ws.Range("E41").Value = "  +14.48%  "
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("E43").Value = "  +10.57%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'17.80"
$numericTextCells += "D44"
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'101.59"
$numericTextCells += "D45"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("D46").Value = "'1.23"
$numericTextCells += "D46"
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("D47").Value = "1.511.86"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D48").Value = "'0.0939"
$numericTextCells += "D48"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "'1.11"
$numericTextCells += "D50"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "'51.58"
$numericTextCells += "D51"
$ws.Range("E51").Value = "  +11.07%  "

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
